# Update "想去人数" (interested-count) figures across sheets, matching the
# regenerated data pulled at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 7612
$ws.Cells.Item(3, 6).Value = 95
$ws.Cells.Item(4, 6).Value = 75
$ws.Cells.Item(5, 6).Value = 5087
$ws.Cells.Item(7, 6).Value = 589
$ws.Cells.Item(8, 6).Value = 611
$ws.Cells.Item(9, 6).Value = 442
$ws.Cells.Item(11, 6).Value = 433
$ws.Cells.Item(12, 6).Value = 759
$ws.Cells.Item(13, 6).Value = 28
$ws.Cells.Item(14, 6).Value = 67
$ws.Cells.Item(15, 6).Value = 261
$ws.Cells.Item(17, 6).Value = 250
$ws.Cells.Item(18, 6).Value = 131
$ws.Cells.Item(19, 6).Value = 383
$ws.Cells.Item(20, 6).Value = 140
$ws.Cells.Item(21, 6).Value = 1081
$ws.Cells.Item(23, 6).Value = 572
$ws.Cells.Item(24, 6).Value = 2163
$ws.Cells.Item(25, 6).Value = 692
$ws.Cells.Item(26, 6).Value = 41
$ws.Cells.Item(27, 6).Value = 39
$ws.Cells.Item(29, 6).Value = 595
$ws.Cells.Item(30, 6).Value = 41

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 313

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 431

# --- Sheet "全部类型" (All types, merged view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 431
$ws.Cells.Item(3, 6).Value = 7612
$ws.Cells.Item(4, 6).Value = 95
$ws.Cells.Item(5, 6).Value = 75
$ws.Cells.Item(7, 6).Value = 5088
$ws.Cells.Item(9, 6).Value = 589
$ws.Cells.Item(10, 6).Value = 611
$ws.Cells.Item(11, 6).Value = 442
$ws.Cells.Item(14, 6).Value = 433
$ws.Cells.Item(15, 6).Value = 313
$ws.Cells.Item(18, 6).Value = 759
$ws.Cells.Item(19, 6).Value = 28
$ws.Cells.Item(20, 6).Value = 67
$ws.Cells.Item(21, 6).Value = 261
$ws.Cells.Item(26, 6).Value = 250
$ws.Cells.Item(27, 6).Value = 131
$ws.Cells.Item(28, 6).Value = 383
$ws.Cells.Item(29, 6).Value = 140
$ws.Cells.Item(30, 6).Value = 1081
$ws.Cells.Item(32, 6).Value = 572
$ws.Cells.Item(33, 6).Value = 2163
$ws.Cells.Item(34, 6).Value = 692
$ws.Cells.Item(35, 6).Value = 41
$ws.Cells.Item(36, 6).Value = 39
$ws.Cells.Item(38, 6).Value = 595
$ws.Cells.Item(39, 6).Value = 41
